$d = $word.ActiveDocument

# Contact-links paragraph currently reads (joining its runs):
#   "portfolio.vercel.app | linkedin.com/in/truongdq | github.com/truong"
# It needs to become:
#   "porfolio-truongdq.vercel.app |  github.com/truong"
# i.e. the portfolio URL changes to "porfolio-truongdq.vercel.app" and the
# "linkedin.com/in/truongdq" segment (plus one of its surrounding spaces)
# is dropped, leaving a double space before "github.com/truong". The
# trailing "nat" run (completing "...truongnat") is untouched.

$old = "portfolio.vercel.app | linkedin.com/in/truongdq | github.com/truong"
$new = "porfolio-truongdq.vercel.app |  github.com/truong"

$range = $d.Content
$found = $range.Find.Execute($old, $true, $false, $false, $false, $false, `
                              $true, 1, $false, $new, 2)

if (-not $found) {
    throw "Could not find target contact-links text to replace"
}

Write-Output "Replace done. Updated text now reads:"
Write-Output $range.Text
